# ---------------------------------------------------------------------------
# Helper: replace a paragraph/run's text with brand-new text while keeping a
# single run (the engine otherwise keeps old/new overlapping substrings as
# separate runs, mirroring a real edit's common-prefix/suffix). Routing the
# assignment through an unrelated placeholder string first guarantees no
# overlap with either the old or the new text, so the result collapses back
# down to one clean <a:r>.
# ---------------------------------------------------------------------------
function Set-CleanText($range, $text) {
    $range.Text = "ZzQx__placeholder__9137"
    $range.Text = $text
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11: "Demonstration - LDD Test Generator"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(3) "Demonstrate a template/mutation-based approach to generating test labels"

# ---------------------------------------------------------------------------
# Slide 13: "Demonstration - Spectral Dictionary Tests"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(7) "You will need to check the output yourself occasionally to verify the overall test suite result."

# ---------------------------------------------------------------------------
# Slide 14: "Interpreting the test output for monolithic tests"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(1) "Monolithic tests require more interpretation that granular tests."
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(2) "The testing system could miss labels that fail for the wrong reason"
$tr = $sh.TextFrame.TextRange
$para3 = $tr.Paragraphs(3)
Set-CleanText $para3 "Log files are built with each push, so check in there for more information"
$para3.IndentLevel = 1

# ---------------------------------------------------------------------------
# Slide 15: "How many tests?"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(1) "There are a lot of things that you can test for your dictionary."
$tr = $sh.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
Set-CleanText $para2 "Some tests are more valuable than others"
$para2.IndentLevel = 1
$tr = $sh.TextFrame.TextRange
$para3 = $tr.Paragraphs(3)
Set-CleanText $para3 "Which tests are the most important to write?"
$para3.IndentLevel = 1

# ---------------------------------------------------------------------------
# Slide 16: "The case against too many tests"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$last = $tr.Paragraphs($tr.Paragraphs().Count)
$last.InsertAfter("`rMore simply, keeping the same structure, but switching between valid values is not a high-value activity.")
$tr = $sh.TextFrame.TextRange
$newPara = $tr.Paragraphs($tr.Paragraphs().Count)
$newPara.IndentLevel = 3

# ---------------------------------------------------------------------------
# Slide 17: "Exercise every class" -> "Testing classes"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(17)
$sh = $s.Shapes.Item(1)
Set-CleanText $sh.TextFrame.TextRange "Testing classes"

$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
# Remove the last paragraph ("Write as many test files...") without leaving
# a stray empty trailing paragraph: copy the (soon to be kept) 3rd
# paragraph's text onto the 4th (currently last) paragraph, then delete the
# no-longer-unique 3rd paragraph -- now a "middle" deletion, which is clean.
$para3 = $tr.Paragraphs(3)
$para4 = $tr.Paragraphs(4)
$para4.Text = $para3.Text
$para3.Delete()

$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(1) "You could write many tests to use each combination of classes, but this is not necessarily valuable."

$tr = $sh.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
Set-CleanText $para2 "It is more valuable to test the minimal description that you can include in each public class."
$tr = $sh.TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
$italicRange = $tr.Characters($para2.Start + 3, 2)
$italicRange.Font.Italic = $true

# ---------------------------------------------------------------------------
# Slide 18: "Exercise every Schematron rule" -> "Exercising Schematron rules"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item(1)
Set-CleanText $sh.TextFrame.TextRange "Exercising Schematron rules"

$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(3).Delete()

$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(1) "An invalid label test could fail on Schematron rules."
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(2) "Valid label tests could pass on Schematron rule, or not trigger the rule at all"
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(3) "Exercising schematron rules is valuable, since they represent exceptions to the exisitng system."

# ---------------------------------------------------------------------------
# Slide 2: "Recap of our goals with testing"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(2) "Ensure that the data that is required is actually captured"
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(3) "Ensure that the schematron rules work"

# ---------------------------------------------------------------------------
# Slide 4: "Principles"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(4) "Minimizing Redundancy"

# ---------------------------------------------------------------------------
# Slide 6: "Monolithic tests vs granular tests"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(1) "Two possible styles of invalid label test are monolithic tests and granular tests"
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(2) "Monolithic tests have few labels with many introduced errors per test"
$tr = $sh.TextFrame.TextRange
$last = $tr.Paragraphs($tr.Paragraphs().Count)
$last.InsertAfter("`rGranular tests have many labels with few introduved errors per test")

# ---------------------------------------------------------------------------
# Slide 7: "Keeping tests granular"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(2) "This simplified interpreting the results, both for you and for the testing framework."

# ---------------------------------------------------------------------------
# Slide 9: "Demonstration - Survey Dictionary Tests"
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
Set-CleanText $tr.Paragraphs(5) "This is enough the trip the validator. When a label fails, it’s for a single reason."
